# Trim the "For more information" section of the abstract: drop the
# "To stay informed about IBM training..." lead-in sentence and the four
# indented social-media lines that followed it (IBM Training News /
# YouTube / Facebook / Twitter), leaving a single empty paragraph with
# a 0.5" left indent (w:ind w:left="720") in their place.

$d = $word.ActiveDocument

# --- Step 1: clear the lead-in sentence's text -----------------------
# Use Find/Replace (ReplaceAll) rather than Range.Delete so the
# paragraph mark (and its pPr: Abstractbodytext / keepNext / keepLines)
# survives intact and becomes the lone remaining paragraph.
$lead = $d.Content
$leadFound = $lead.Find.Execute(
    "To stay informed about IBM training, see the following sites:",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# --- Step 2: remove the four social-media paragraphs entirely --------
# Locate the exact start of the first line and the exact end of the
# last line by text match (robust to any offset shifts from step 1),
# then delete the whole span, including every paragraph mark in it, so
# the block collapses away completely.
if ($leadFound) {
    $startR = $d.Content
    $startFound = $startR.Find.Execute(
        "IBM Training News: https://www.ibm.com/blogs/ibm-training",
        $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

    $endR = $d.Content
    $endFound = $endR.Find.Execute(
        "Twitter: twitter.com/ibm",
        $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

    if ($startFound -and $endFound) {
        $blockStart = $startR.Start
        $blockEnd = $endR.End + 1   # +1 swallows the final paragraph mark

        $victim = $d.Range($blockStart, $blockEnd)
        $victim.Delete()
    }

    # --- Step 3: give the surviving (now empty) paragraph the same
    # 0.5" left indent (36 pt == w:ind w:left="720") the removed
    # indented lines had.
    $survivor = $d.Paragraphs.Item($d.Paragraphs.Count)
    $survivor.LeftIndent = 36
}
